$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A, rows 105 through 143, from 89.29531013618748 to 63.80071144077566
$ws.Range("A105:A143").Value = 63.80071144077566
